$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "Weights: [-1/(sqrt(nc); -1/(sqrt(nc)] where nc is the number of inputs
# for" : the runs "-1/(sqrt(nc)" and "] where nc is the number of inputs
# for " get merged into a single run. A same-text single Find/Replace over
# that span (without touching the separate leading-space run before it)
# reproduces exactly that run merge.
$ok1 = $d.Content.Find.Execute("where nc is", $true, $false, $false, $false, $false, $true, 1, $false, "where nc is", 1)

# --- Change 2 -------------------------------------------------------------
# Fix the "vakpropagation" typo to "back propagation" and move the
# "_GoBack" bookmark so that it now sits between "...does the back " and
# "propagation...".
$okTypo = $d.Content.Find.Execute("vak", $true, $false, $false, $false, $false, $true, 1, $false, "back ", 1)

# Split the run right after "computed activ" (first new run boundary) using
# a throw-away bookmark collapsed at that point - Bookmarks.Add on a
# collapsed range splits the run it lands in without altering any text.
$text = $d.Content.Text
$marker = "bprop will use computed activ"
$idx1 = $text.IndexOf($marker)
$split1 = $idx1 + $marker.Length
$r1 = $d.Range($split1, $split1)
$d.Bookmarks.Add("ZZZTempSplit", $r1)

# Now place "_GoBack" at the second split point, right after
# "...does the back ". Since bookmark names are unique, adding it here
# also removes it from its old location further down the document.
$text2 = $d.Content.Text
$marker2 = "bprop will use computed activations by fprop and does the back "
$idx2 = $text2.IndexOf($marker2)
$split2 = $idx2 + $marker2.Length
$r2 = $d.Range($split2, $split2)
$d.Bookmarks.Add("_GoBack", $r2)

# Remove the temporary bookmark - the run split it created stays in place.
$tmp = $d.Bookmarks.Item("ZZZTempSplit")
$tmp.Delete()
